# The commit swaps the OOXML content of ppt/theme/theme1.xml (the
# default "Office Theme" used by the Notes Master) and
# ppt/theme/theme2.xml (the "Integral" theme used by the Slide Master /
# actual deck) so that:
#   - theme1.xml ends up holding the "Integral" colour scheme
#   - theme2.xml ends up holding the default "Office Theme" colour scheme
#
# The font scheme (<a:fontScheme>) and format scheme (<a:fmtScheme>)
# blocks are byte-for-byte identical between the two themes already, so
# the only real content difference is the <a:clrScheme> (the 12 theme
# colours) plus the cosmetic name="" attributes on <a:theme>/<a:clrScheme>.
#
# The PowerPoint object model only exposes one live "theme" to script -
# the one driving the slide master / slides (backed by theme2.xml); the
# Notes Master's theme part isn't independently reachable through the
# COM surface. So we apply the colour half of the swap - which is the
# part with visible, functional effect - to that theme via
# ThemeColorScheme, setting every slot to the current Office Theme
# default values (i.e. what theme2.xml's colours should become).

$p = $ppt.ActivePresentation

function ColorVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Default "Office Theme" colours (current ppt/theme/theme1.xml), in
# ThemeColorScheme slot order.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = ColorVal $officeThemeColors[$i - 1]
}
